$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value2 = 4553.8335  # H76
$ws.Cells.Item(76, 9).Value2 = 3806.3333  # I76
$ws.Cells.Item(76, 10).Value2 = 5301.3335  # J76
$ws.Cells.Item(76, 11).Value2 = 3806.3333  # K76
$ws.Cells.Item(76, 12).Value2 = 5301.3335  # L76
$ws.Cells.Item(76, 13).Value2 = -3491.3333  # M76
$ws.Cells.Item(76, 14).Value2 = -5931.3335  # N76
$ws.Cells.Item(79, 8).Value2 = 4553.8335  # H79
$ws.Cells.Item(79, 9).Value2 = 3806.3333  # I79
$ws.Cells.Item(79, 10).Value2 = 5301.3335  # J79
$ws.Cells.Item(79, 11).Value2 = 3806.3333  # K79
$ws.Cells.Item(79, 12).Value2 = 5301.3335  # L79
$ws.Cells.Item(79, 13).Value2 = -2714.3333  # M79
$ws.Cells.Item(79, 14).Value2 = -7485.3335  # N79
$ws.Cells.Item(86, 8).Value2 = 1976.4445  # H86
$ws.Cells.Item(86, 9).Value2 = 1782.8948  # I86
$ws.Cells.Item(86, 11).Value2 = 1782.8948  # K86
$ws.Cells.Item(86, 13).Value2 = -659.8948  # M86
$ws.Cells.Item(89, 8).Value2 = 1976.4445  # H89
$ws.Cells.Item(89, 9).Value2 = 1782.8948  # I89
$ws.Cells.Item(89, 11).Value2 = 8914.474  # K89
$ws.Cells.Item(89, 13).Value2 = -3298.474  # M89
$ws.Cells.Item(106, 8).Value2 = 2926.625  # H106
$ws.Cells.Item(106, 9).Value2 = 2285.6  # I106
$ws.Cells.Item(106, 11).Value2 = 2285.6  # K106
$ws.Cells.Item(106, 13).Value2 = -1654.6  # M106
$ws.Cells.Item(112, 8).Value2 = 8473  # H112
$ws.Cells.Item(112, 10).Value2 = 8858.387  # J112
$ws.Cells.Item(112, 12).Value2 = 26575.161  # L112
$ws.Cells.Item(112, 14).Value2 = -28791.161  # N112
$ws.Cells.Item(133, 9).Value2 = 0  # I133
$ws.Cells.Item(133, 10).Value2 = 107498.25  # J133
$ws.Cells.Item(133, 11).Value2 = 0  # K133
$ws.Cells.Item(133, 12).Value2 = 107498.25  # L133
$ws.Cells.Item(133, 13).Value = $null  # M133
$ws.Cells.Item(133, 14).Value2 = -117618.25  # N133
$ws.Cells.Item(137, 8).Value2 = 11408.765  # H137
$ws.Cells.Item(137, 9).Value2 = 19608.646  # I137
$ws.Cells.Item(137, 10).Value2 = 3208.8823  # J137
$ws.Cells.Item(137, 11).Value2 = 58825.938  # K137
$ws.Cells.Item(137, 12).Value2 = 9626.6469  # L137
$ws.Cells.Item(137, 13).Value2 = -56275.938  # M137
$ws.Cells.Item(137, 14).Value2 = -14726.6469  # N137
$ws.Cells.Item(138, 8).Value2 = 2847.359  # H138
$ws.Cells.Item(138, 9).Value2 = 1851.2632  # I138
$ws.Cells.Item(138, 11).Value2 = 5553.7896  # K138
$ws.Cells.Item(138, 13).Value2 = -413.7896000000001  # M138

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 21446.178  # H32
$ws.Cells.Item(32, 10).Value2 = 6431.5  # J32
$ws.Cells.Item(32, 12).Value2 = 6431.5  # L32
$ws.Cells.Item(32, 14).Value2 = -7005.5  # N32
$ws.Cells.Item(132, 8).Value2 = 1266.0444  # H132
$ws.Cells.Item(132, 9).Value2 = 884  # I132
$ws.Cells.Item(132, 11).Value2 = 2652  # K132
$ws.Cells.Item(132, 13).Value2 = -122  # M132
$ws.Cells.Item(139, 8).Value2 = 81000  # H139
$ws.Cells.Item(139, 10).Value2 = 81000  # J139
$ws.Cells.Item(139, 12).Value2 = 81000  # L139
$ws.Cells.Item(139, 14).Value2 = -91280  # N139

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(75, 8).Value2 = 30110.5  # H75
$ws.Cells.Item(75, 9).Value2 = 6000  # I75
$ws.Cells.Item(75, 11).Value2 = 6000  # K75
$ws.Cells.Item(75, 13).Value2 = -5064  # M75
$ws.Cells.Item(78, 8).Value2 = 30110.5  # H78
$ws.Cells.Item(78, 9).Value2 = 6000  # I78
$ws.Cells.Item(78, 11).Value2 = 18000  # K78
$ws.Cells.Item(78, 13).Value2 = -13320  # M78
$ws.Cells.Item(105, 8).Value2 = 1872.4  # H105
$ws.Cells.Item(105, 9).Value2 = 1125.5625  # I105
$ws.Cells.Item(105, 11).Value2 = 1125.5625  # K105
$ws.Cells.Item(105, 13).Value2 = 621.4375  # M105
$ws.Cells.Item(107, 8).Value2 = 54268.2  # H107
$ws.Cells.Item(107, 9).Value2 = 66585.375  # I107
$ws.Cells.Item(107, 10).Value2 = 4999.5  # J107
$ws.Cells.Item(107, 11).Value2 = 66585.375  # K107
$ws.Cells.Item(107, 12).Value2 = 4999.5  # L107
$ws.Cells.Item(107, 13).Value2 = -64665.375  # M107
$ws.Cells.Item(107, 14).Value2 = -8839.5  # N107

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value2 = 1575  # H16
$ws.Cells.Item(16, 9).Value2 = 1168.5714  # I16
$ws.Cells.Item(16, 11).Value2 = 1168.5714  # K16
$ws.Cells.Item(16, 13).Value2 = -881.5714  # M16
$ws.Cells.Item(31, 8).Value2 = 2633878.2  # H31
$ws.Cells.Item(31, 9).Value2 = 3228043.2  # I31
$ws.Cells.Item(31, 10).Value2 = 2576.1428  # J31
$ws.Cells.Item(31, 11).Value2 = 3228043.2  # K31
$ws.Cells.Item(31, 12).Value2 = 2576.1428  # L31
$ws.Cells.Item(31, 13).Value2 = -3227748.2  # M31
$ws.Cells.Item(31, 14).Value2 = -3166.1428  # N31
$ws.Cells.Item(34, 8).Value2 = 2633878.2  # H34
$ws.Cells.Item(34, 9).Value2 = 3228043.2  # I34
$ws.Cells.Item(34, 10).Value2 = 2576.1428  # J34
$ws.Cells.Item(34, 11).Value2 = 3228043.2  # K34
$ws.Cells.Item(34, 12).Value2 = 2576.1428  # L34
$ws.Cells.Item(34, 13).Value2 = -3227841.2  # M34
$ws.Cells.Item(34, 14).Value2 = -2980.1428  # N34
$ws.Cells.Item(50, 8).Value2 = 36999.668  # H50
$ws.Cells.Item(50, 10).Value2 = 36999.668  # J50
$ws.Cells.Item(50, 12).Value2 = 36999.668  # L50
$ws.Cells.Item(50, 14).Value2 = -38249.668  # N50
$ws.Cells.Item(51, 8).Value2 = 30659.6  # H51
$ws.Cells.Item(51, 10).Value2 = 30659.6  # J51
$ws.Cells.Item(51, 12).Value2 = 30659.6  # L51
$ws.Cells.Item(51, 14).Value2 = -32131.6  # N51
$ws.Cells.Item(60, 8).Value2 = 34732.918  # H60
$ws.Cells.Item(60, 10).Value2 = 34779.8  # J60
$ws.Cells.Item(60, 12).Value2 = 34779.8  # L60
$ws.Cells.Item(60, 14).Value2 = -35801.8  # N60
$ws.Cells.Item(61, 8).Value2 = 30659.6  # H61
$ws.Cells.Item(61, 10).Value2 = 30659.6  # J61
$ws.Cells.Item(61, 12).Value2 = 30659.6  # L61
$ws.Cells.Item(61, 14).Value2 = -31355.6  # N61
$ws.Cells.Item(107, 8).Value2 = 648.34485  # H107
$ws.Cells.Item(107, 9).Value2 = 483.6316  # I107
$ws.Cells.Item(107, 10).Value2 = 961.3  # J107
$ws.Cells.Item(107, 11).Value2 = 483.6316  # K107
$ws.Cells.Item(107, 12).Value2 = 961.3  # L107
$ws.Cells.Item(107, 13).Value2 = 1436.3684  # M107
$ws.Cells.Item(107, 14).Value2 = -4801.3  # N107
$ws.Cells.Item(113, 8).Value2 = 1575  # H113
$ws.Cells.Item(113, 9).Value2 = 1168.5714  # I113
$ws.Cells.Item(113, 11).Value2 = 1168.5714  # K113
$ws.Cells.Item(113, 13).Value2 = 1001.4286  # M113
$ws.Cells.Item(132, 8).Value2 = 20173.281  # H132
$ws.Cells.Item(132, 9).Value2 = 22896.572  # I132
$ws.Cells.Item(132, 10).Value2 = 1110.25  # J132
$ws.Cells.Item(132, 11).Value2 = 68689.716  # K132
$ws.Cells.Item(132, 12).Value2 = 3330.75  # L132
$ws.Cells.Item(132, 13).Value2 = -66159.716  # M132
$ws.Cells.Item(132, 14).Value2 = -8390.75  # N132
$ws.Cells.Item(134, 8).Value2 = 2558.1738  # H134
$ws.Cells.Item(134, 9).Value2 = 2492.6365  # I134
$ws.Cells.Item(134, 10).Value2 = 4000  # J134
$ws.Cells.Item(134, 11).Value2 = 7477.9095  # K134
$ws.Cells.Item(134, 12).Value2 = 12000  # L134
$ws.Cells.Item(134, 13).Value2 = -4942.9095  # M134
$ws.Cells.Item(134, 14).Value2 = -17070  # N134
$ws.Cells.Item(138, 8).Value2 = 112301.71  # H138
$ws.Cells.Item(138, 10).Value2 = 112301.71  # J138
$ws.Cells.Item(138, 12).Value2 = 112301.71  # L138
$ws.Cells.Item(138, 14).Value2 = -122581.71  # N138

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value2 = 1679.4  # H132
$ws.Cells.Item(132, 9).Value2 = 1934.3334  # I132
$ws.Cells.Item(132, 11).Value2 = 17409.0006  # K132
$ws.Cells.Item(132, 13).Value2 = -14879.0006  # M132
$ws.Cells.Item(136, 8).Value2 = 5248.067  # H136
$ws.Cells.Item(136, 9).Value2 = 2391.2222  # I136
$ws.Cells.Item(136, 11).Value2 = 7173.6666  # K136
$ws.Cells.Item(136, 13).Value2 = -2073.6666  # M136
$ws.Cells.Item(138, 8).Value2 = 10674.6875  # H138
$ws.Cells.Item(138, 9).Value2 = 11080.2  # I138
$ws.Cells.Item(138, 10).Value2 = 9998.833  # J138
$ws.Cells.Item(138, 11).Value2 = 33240.60000000001  # K138
$ws.Cells.Item(138, 12).Value2 = 29996.499  # L138
$ws.Cells.Item(138, 13).Value2 = -28100.60000000001  # M138
$ws.Cells.Item(138, 14).Value2 = -40276.499  # N138
$ws.Cells.Item(141, 8).Value2 = 5792.3335  # H141
$ws.Cells.Item(141, 9).Value2 = 5934.75  # I141
$ws.Cells.Item(141, 11).Value2 = 17804.25  # K141
$ws.Cells.Item(141, 13).Value2 = -12624.25  # M141

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value2 = 155.71428  # H2
$ws.Cells.Item(2, 9).Value2 = 55.8  # I2
$ws.Cells.Item(2, 11).Value2 = 55.8  # K2
$ws.Cells.Item(2, 13).Value2 = 57.2  # M2
$ws.Cells.Item(44, 8).Value2 = 5022583.5  # H44
$ws.Cells.Item(44, 10).Value2 = 7520000  # J44
$ws.Cells.Item(44, 12).Value2 = 7520000  # L44
$ws.Cells.Item(44, 14).Value2 = -7521192  # N44
$ws.Cells.Item(126, 8).Value2 = 2342.0625  # H126
$ws.Cells.Item(126, 9).Value2 = 2092.818  # I126
$ws.Cells.Item(126, 10).Value2 = 2890.4  # J126
$ws.Cells.Item(126, 11).Value2 = 6278.454000000001  # K126
$ws.Cells.Item(126, 12).Value2 = 8671.2  # L126
$ws.Cells.Item(126, 13).Value2 = -3808.454000000001  # M126
$ws.Cells.Item(126, 14).Value2 = -13611.2  # N126
$ws.Cells.Item(132, 8).Value2 = 5113.533  # H132
$ws.Cells.Item(132, 9).Value2 = 5220.8965  # I132
$ws.Cells.Item(132, 11).Value2 = 15662.6895  # K132
$ws.Cells.Item(132, 13).Value2 = -13132.6895  # M132
$ws.Cells.Item(139, 8).Value2 = 104993.336  # H139
$ws.Cells.Item(139, 10).Value2 = 104993.336  # J139
$ws.Cells.Item(139, 12).Value2 = 104993.336  # L139
$ws.Cells.Item(139, 14).Value2 = -115273.336  # N139
$ws.Cells.Item(140, 8).Value2 = 110253  # H140
$ws.Cells.Item(140, 10).Value2 = 110253  # J140
$ws.Cells.Item(140, 12).Value2 = 110253  # L140
$ws.Cells.Item(140, 14).Value2 = -120613  # N140

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value2 = 2745.5625  # H7
$ws.Cells.Item(7, 9).Value2 = 2552.1428  # I7
$ws.Cells.Item(7, 10).Value2 = 4099.5  # J7
$ws.Cells.Item(7, 11).Value2 = 2552.1428  # K7
$ws.Cells.Item(7, 12).Value2 = 4099.5  # L7
$ws.Cells.Item(7, 13).Value2 = -2440.1428  # M7
$ws.Cells.Item(7, 14).Value2 = -4323.5  # N7
$ws.Cells.Item(115, 8).Value2 = 50000  # H115
$ws.Cells.Item(115, 10).Value2 = 50000  # J115
$ws.Cells.Item(115, 12).Value2 = 50000  # L115
$ws.Cells.Item(115, 14).Value2 = -52350  # N115
$ws.Cells.Item(126, 8).Value2 = 2745.5625  # H126
$ws.Cells.Item(126, 9).Value2 = 2552.1428  # I126
$ws.Cells.Item(126, 10).Value2 = 4099.5  # J126
$ws.Cells.Item(126, 11).Value2 = 7656.428400000001  # K126
$ws.Cells.Item(126, 12).Value2 = 12298.5  # L126
$ws.Cells.Item(126, 13).Value2 = -5186.428400000001  # M126
$ws.Cells.Item(126, 14).Value2 = -17238.5  # N126

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value2 = 37870.625  # H2
$ws.Cells.Item(2, 9).Value2 = 45050  # I2
$ws.Cells.Item(2, 11).Value2 = 45050  # K2
$ws.Cells.Item(2, 13).Value2 = -44938  # M2
$ws.Cells.Item(4, 8).Value2 = 770916.6  # H4
$ws.Cells.Item(4, 9).Value2 = 1277.6666  # I4
$ws.Cells.Item(4, 11).Value2 = 1277.6666  # K4
$ws.Cells.Item(4, 13).Value2 = -1164.6666  # M4
$ws.Cells.Item(11, 8).Value2 = 100000  # H11
$ws.Cells.Item(11, 9).Value2 = 0  # I11
$ws.Cells.Item(11, 10).Value2 = 100000  # J11
$ws.Cells.Item(11, 11).Value2 = 0  # K11
$ws.Cells.Item(11, 12).Value2 = 100000  # L11
$ws.Cells.Item(11, 13).Value = $null  # M11
$ws.Cells.Item(11, 14).Value2 = -100284  # N11
$ws.Cells.Item(132, 8).Value2 = 3446.0293  # H132
$ws.Cells.Item(132, 9).Value2 = 3392.3667  # I132
$ws.Cells.Item(132, 11).Value2 = 10177.1001  # K132
$ws.Cells.Item(132, 13).Value2 = -7647.1001  # M132
